$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Simple single-cell corrections ---
$ws.Range("C2").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("D15").Value = "-"

# --- Row 18 ---
$ws.Range("B18").Value = "[Weslei-Metrologia-1NA, -, -, -]"
$ws.Range("C18").Value = "[Cleidson-Sistemas digitais-1NA, Pedro Bispo-Manut. Elétrica-1NA, Sandro-Lógica de Programação-1NA, Sandro-Acionamentos Elétricos-1NA]"
$ws.Range("D18").Value = "[Suzanny-Des. Bas. Mec.-1NA, Suzanny-Des. Bas. Mec.-1NA]"
$ws.Range("E18").Value = "[Weslei-Metrologia-1NA, Rogério-Processos de Usinagem 2-1NA, Anderson-Processos de Usinagem 1-1NA, Leonardo-Manut. Mecânica-1NA]"
$ws.Range("F18").Value = "[Rachel-Tecnologia dos Materiais.-1NA, Rachel-Tecnologia dos Materiais.-1NA]"

# --- Row 19 ---
$ws.Range("B19").Value = "[Weslei-Metrologia-1NA, -, -, -]"
$ws.Range("C19").Value = "[Cleidson-Sistemas digitais-1NA, Pedro Bispo-Manut. Elétrica-1NA, Sandro-Lógica de Programação-1NA, Sandro-Acionamentos Elétricos-1NA]"
$ws.Range("D19").Value = "Allan Cupertino-Circuitos Elétricos 1-"
$ws.Range("E19").Value = "[Weslei-Metrologia-1NA, Rogério-Processos de Usinagem 2-1NA, Anderson-Processos de Usinagem 1-1NA, Leonardo-Manut. Mecânica-1NA]"
$ws.Range("F19").Value = "[Rachel-Tecnologia dos Materiais.-1NA, Rachel-Tecnologia dos Materiais.-1NA]"

# --- Row 20 ---
$ws.Range("B20").Value = "Euclides-Gestão Integrada-"
$ws.Range("C20").Value = "[Cleidson-Sistemas digitais-1NA, Pedro Bispo-Manut. Elétrica-1NA, Sandro-Lógica de Programação-1NA, Sandro-Acionamentos Elétricos-1NA]"
# D20 is unchanged: [Suzanny-Des. Bas. Mec.-1NA, Suzanny-Des. Bas. Mec.-1NA]
$ws.Range("E20").Value = "Allan Cupertino-Circuitos Elétricos 1-"
$ws.Range("F20").Value = "Euclides-Gestão Integrada-"

# --- Row 21 ---
$ws.Range("B21").Value = "Gilberto-Tecnologias Mecânicas-"
$ws.Range("C21").Value = "[Cleidson-Sistemas digitais-1NA, Pedro Bispo-Manut. Elétrica-1NA, Sandro-Lógica de Programação-1NA, Sandro-Acionamentos Elétricos-1NA]"
$ws.Range("D21").Value = "[Suzanny-Des. Bas. Mec.-1NA, Suzanny-Des. Bas. Mec.-1NA]"
$ws.Range("E21").Value = "[-, Rogério-Processos de Usinagem 2-1NA, Anderson-Processos de Usinagem 1-1NA, Leonardo-Manut. Mecânica-1NA]"
$ws.Range("F21").Value = "Gilberto-Tecnologias Mecânicas-"
